$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data edits on row 3 (trial 2): x_corrSteps, y_nrSteps, alienID
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 13

# Move the active selection from A6 to E3, matching the saved cursor position
$ws.Range("E3").Select()
